# Expenses report: add the "Employee Hours" mini-table below the existing
# product-sales section, and a total-cost formula under it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 39): Date | Employee | Role | Start Shift | End Shift | Total Hours | Salary/Hour | Cost
$ws.Range("B39").Value = "Date"
$ws.Range("C39").Value = "Employee"
$ws.Range("D39").Value = "Role"
$ws.Range("E39").Value = "Start Shift"
$ws.Range("F39").Value = "End Shift"
$ws.Range("G39").Value = "Total Hours"
$ws.Range("H39").Value = "Salary/Hour"
$ws.Range("I39").Value = "Cost"

# --- Start/End Shift columns are meant to hold clock times (row 40 onward);
# apply the HH:MM time format used for shift entry, then clear the scratch
# cells back out since no shift rows have been entered yet.
$ws.Range("E40").NumberFormat = "HH:MM"
$ws.Range("F40").NumberFormat = "HH:MM"
$ws.Range("E40:F40").Clear()

# --- Total cost formula (row 41), formatted like the other currency totals.
$ws.Range("I41").Formula = "=SUM(I41:I41)"
$ws.Range("I41").NumberFormat = $ws.Range("J5").NumberFormat
